# "admin authenticated services updated"
# Update the vouchers table: names (B), voucher type (D), price as a real
# number instead of a text placeholder (E), and expiry dates (F).
# Column A (id) and C (code) are unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names   = @("Burger King", "Zamber", "Renee", "Yatra", "Ibibo", "Sugar")
$types   = @("dining", "dining", "cosmetics", "travel", "travel", "cosmetics")
$prices  = @(250, 300, 269, 200, 199, 299)
$dates   = @(45178, 45178, 45179, 45177, 45178, 45177)

for ($i = 0; $i -lt 6; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 4).Value = $types[$i]
    $ws.Cells.Item($row, 5).Value = $prices[$i]
    $ws.Cells.Item($row, 6).Value = $dates[$i]
}

$ws.Range("I8").Select()
